# Update Name of Algo
# Apply updated imputed values to result_data_KNN sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.591
$ws.Range("C3").Value = -12.553
$ws.Range("C5").Value = -12.261
$ws.Range("E7").Value = 12.752
$ws.Range("B9").Value = 6.886
$ws.Range("E9").Value = 12.787
$ws.Range("C11").Value = -12.917
$ws.Range("C12").Value = -12.628
$ws.Range("B13").Value = 6.308000000000001
$ws.Range("B16").Value = 5.77
$ws.Range("B18").Value = 5.91
$ws.Range("B20").Value = 6.308000000000001
$ws.Range("C21").Value = -12.747
$ws.Range("E21").Value = 12.925
